$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Prepare formatting for the new rows before writing values -----------
# Column A (dates) for rows 29-34 should match the existing date format
# used by the rows above (e.g. A28).
$ws.Range("A28").Copy()
$ws.Range("A29:A34").PasteSpecial(-4122)  # xlPasteFormats

# Columns B:F for rows 29-38 should match the blank "template" row format
# (e.g. row 2), which is style 1/1/14/14/1/1 for A/B/C/D/E/F.
$ws.Range("B2:F2").Copy()
$ws.Range("B29:F38").PasteSpecial(-4122)  # xlPasteFormats

# Column A for the trailing blank rows 35-38 should match the blank
# template row's own (non-date) style as well.
$ws.Range("A2").Copy()
$ws.Range("A35:A38").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Row 29 ---------------------------------------------------------------
$ws.Range("A29").Value = 45743
$ws.Range("B29").Value = "Odata"
$ws.Range("C29").Value = "Minor changes in Pit stop"
$ws.Range("D29").Value = "HAS IOT , BATTSLNO, Battery model and other small changes"
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = "Don't know if we need to capture the efforts "

# --- Row 30 ---------------------------------------------------------------
$ws.Range("A30").Value = 45744
$ws.Range("B30").Value = "All master"
$ws.Range("C30").Value = "Emergency corrections"
$ws.Range("D30").Value = "Core class, FM's etc has been changed by Roopesh, so fixing all that"
$ws.Range("E30").Value = 3

# --- Row 31 ---------------------------------------------------------------
$ws.Range("A31").Value = 45745
$ws.Range("D31").Value = "Not much progress"

# --- Row 32 ---------------------------------------------------------------
$ws.Range("A32").Value = 45746
$ws.Range("D32").Value = "Not much progress"

# --- Row 33 ---------------------------------------------------------------
$ws.Range("A33").Value = 45747
$ws.Range("C33").Value = "Naming convention corrections  and BOM Master"
$ws.Range("D33").Value = "Naming convention changes and BOM Master develoopment "
$ws.Range("E33").Value = "2 - 3 hrs"
$ws.Range("F33").Value = "Corrected in 6 masters  - Check if we need to capture this efforts"

# --- Row 34 ---------------------------------------------------------------
$ws.Range("A34").Value = 45748
$ws.Range("B34").Value = "Odata "
$ws.Range("C34").Value = "Odata changes"
$ws.Range("D34").Value = "PitStop master, Vehicle allotment changes"
$ws.Range("E34").Value = 2

# --- View state: active cell / selection ----------------------------------
$ws.Activate()
$ws.Range("D33").Select()
